# Auto-generated edit script: applies numeric cell updates across all 8 leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (65 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3875.25
$ws.Range("I2").Value = 4667
$ws.Range("K2").Value = 4667
$ws.Range("M2").Value = -4554
$ws.Range("H32").Value = 7070.7144
$ws.Range("J32").Value = 8099
$ws.Range("L32").Value = 8099
$ws.Range("N32").Value = -8751
$ws.Range("H39").Value = 2275
$ws.Range("I39").Value = 279.64285
$ws.Range("J39").Value = 9258.75
$ws.Range("K39").Value = 838.9285500000001
$ws.Range("L39").Value = 27776.25
$ws.Range("M39").Value = -542.9285500000001
$ws.Range("N39").Value = -28368.25
$ws.Range("H41").Value = 315.08334
$ws.Range("I41").Value = 397.57144
$ws.Range("J41").Value = 199.6
$ws.Range("K41").Value = 397.57144
$ws.Range("L41").Value = 199.6
$ws.Range("M41").Value = 42.42856
$ws.Range("N41").Value = -1079.6
$ws.Range("H86").Value = 6800
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877
$ws.Range("H89").Value = 6800
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H98").Value = 3499
$ws.Range("I98").Value = 2665
$ws.Range("K98").Value = 2665
$ws.Range("M98").Value = -1167
$ws.Range("H112").Value = 2436.95
$ws.Range("J112").Value = 2436.95
$ws.Range("L112").Value = 7310.849999999999
$ws.Range("N112").Value = -9526.849999999999
$ws.Range("H122").Value = 3499
$ws.Range("I122").Value = 2665
$ws.Range("K122").Value = 7995
$ws.Range("M122").Value = -5545
$ws.Range("H132").Value = 1071.151
$ws.Range("J132").Value = 2469.7144
$ws.Range("L132").Value = 7409.1432
$ws.Range("N132").Value = -12469.1432
$ws.Range("H135").Value = 1208.659
$ws.Range("I135").Value = 587.6061
$ws.Range("K135").Value = 5288.4549
$ws.Range("M135").Value = -2753.4549
$ws.Range("H138").Value = 5231.4653
$ws.Range("I138").Value = 4541.231
$ws.Range("J138").Value = 5430.8667
$ws.Range("K138").Value = 13623.693
$ws.Range("L138").Value = 16292.6001
$ws.Range("M138").Value = -8483.692999999999
$ws.Range("N138").Value = -26572.6001
$ws.Range("H141").Value = 2493.7
$ws.Range("I141").Value = 2215.2222
$ws.Range("K141").Value = 6645.6666
$ws.Range("M141").Value = -1465.6666

# --- Sheet: ARM (19 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20141.264
$ws.Range("I32").Value = 20131.986
$ws.Range("K32").Value = 20131.986
$ws.Range("M32").Value = -19844.986
$ws.Range("H74").Value = 273342.53
$ws.Range("I74").Value = 295402.22
$ws.Range("K74").Value = 295402.22
$ws.Range("M74").Value = -294528.22
$ws.Range("H77").Value = 273342.53
$ws.Range("I77").Value = 295402.22
$ws.Range("K77").Value = 1477011.1
$ws.Range("M77").Value = -1472643.1
$ws.Range("H97").Value = 3368555.2
$ws.Range("I97").Value = 3368555.2
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3368555.2
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -3368059.2
$ws.Range("N97").ClearContents()

# --- Sheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2057.2856
$ws.Range("I36").Value = 2057.2856
$ws.Range("K36").Value = 2057.2856
$ws.Range("M36").Value = -1523.2856
$ws.Range("H94").Value = 2158.7693
$ws.Range("I94").Value = 1128.8
$ws.Range("K94").Value = 1128.8
$ws.Range("M94").Value = -677.8

# --- Sheet: CRP (16 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 168189.67
$ws.Range("I12").Value = 250050
$ws.Range("K12").Value = 250050
$ws.Range("M12").Value = -249880
$ws.Range("H31").Value = 22228032
$ws.Range("I31").Value = 52634710
$ws.Range("K31").Value = 52634710
$ws.Range("M31").Value = -52634415
$ws.Range("H34").Value = 22228032
$ws.Range("I34").Value = 52634710
$ws.Range("K34").Value = 52634710
$ws.Range("M34").Value = -52634508
$ws.Range("H141").Value = 191830.23
$ws.Range("J141").Value = 220796.31
$ws.Range("L141").Value = 220796.31
$ws.Range("N141").Value = -231156.31

# --- Sheet: CUL (27 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 24376562
$ws.Range("I4").Value = 34385776
$ws.Range("K4").Value = 103157328
$ws.Range("M4").Value = -103157216
$ws.Range("H7").Value = 450
$ws.Range("I7").Value = 466.66666
$ws.Range("K7").Value = 1399.99998
$ws.Range("M7").Value = -1287.99998
$ws.Range("H36").Value = 608.1111
$ws.Range("J36").Value = 415
$ws.Range("L36").Value = 1245
$ws.Range("N36").Value = -1583
$ws.Range("H57").Value = 3916.2
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 3916.2
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 11748.6
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -12866.6
$ws.Range("H81").Value = 500
$ws.Range("I81").Value = 500
$ws.Range("K81").Value = 1500
$ws.Range("M81").Value = -377
$ws.Range("H84").Value = 500
$ws.Range("I84").Value = 500
$ws.Range("K84").Value = 4500
$ws.Range("M84").Value = 1116

# --- Sheet: GSM (18 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1568.9
$ws.Range("I97").Value = 1458.6
$ws.Range("K97").Value = 1458.6
$ws.Range("M97").Value = -962.5999999999999
$ws.Range("H102").Value = 2389.5
$ws.Range("I102").Value = 1601.375
$ws.Range("J102").Value = 3965.75
$ws.Range("K102").Value = 1601.375
$ws.Range("L102").Value = 3965.75
$ws.Range("M102").Value = 20.625
$ws.Range("N102").Value = -7209.75
$ws.Range("H113").Value = 4129.2
$ws.Range("I113").Value = 3786.5
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 3786.5
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -1616.5
$ws.Range("N113").Value = -9840

# --- Sheet: LTW (32 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 125004000
$ws.Range("I40").Value = 125004000
$ws.Range("K40").Value = 125004000
$ws.Range("M40").Value = -125003864
$ws.Range("H61").Value = 4754.25
$ws.Range("I61").Value = 4754.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4754.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4552.25
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 4754.25
$ws.Range("I113").Value = 4754.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4754.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2584.25
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 45552530
$ws.Range("I122").Value = 62505724
$ws.Range("J122").Value = 344000
$ws.Range("K122").Value = 187517172
$ws.Range("L122").Value = 1032000
$ws.Range("M122").Value = -187514722
$ws.Range("N122").Value = -1036900
$ws.Range("H136").Value = 5861.037
$ws.Range("I136").Value = 4335.923
$ws.Range("J136").Value = 9826.333000000001
$ws.Range("K136").Value = 13007.769
$ws.Range("L136").Value = 29478.999
$ws.Range("M136").Value = -10457.769
$ws.Range("N136").Value = -34578.999

# --- Sheet: WVR (27 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29500
$ws.Range("J64").Value = 29500
$ws.Range("L64").Value = 29500
$ws.Range("N64").Value = -29996
$ws.Range("H67").Value = 29500
$ws.Range("J67").Value = 29500
$ws.Range("L67").Value = 29500
$ws.Range("N67").Value = -31216
$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678
$ws.Range("H122").Value = 2670.434
$ws.Range("I122").Value = 2690.861
$ws.Range("K122").Value = 8072.583
$ws.Range("M122").Value = -5622.583
$ws.Range("H126").Value = 3649.0715
$ws.Range("I126").Value = 3829.7693
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 11489.3079
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -9019.3079
$ws.Range("N126").Value = -8840
$ws.Range("H136").Value = 3391.5
$ws.Range("I136").Value = 2274.0454
$ws.Range("K136").Value = 6822.1362
$ws.Range("M136").Value = -4272.1362

